$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for "runs, balls, fours, sixes" (columns C:F) got reshuffled
# between rows: row2<->row5, row4<->row8, row6<->row7 (row3 unchanged).
# Force a text number format before writing so the values stay text,
# matching the original numberStoredAsText layout.
$ws.Range("C2:F8").NumberFormat = "@"

function Set-Row($row, $c, $d, $e, $f) {
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
}

Set-Row 2 "66" "56" "8" "2"
Set-Row 4 "6"  "6"  "0" "0"
Set-Row 5 "0"  "1"  "0" "0"
Set-Row 6 "0"  "1"  "0" "0"
Set-Row 7 "27" "16" "0" "2"
Set-Row 8 "17" "14" "1" "0"
